# Scenario_100: add a new "Gap Naive" column to the TestingFile results
# table, reorder "Iterations"/"Converged?" ahead of "Gap LR", and correct a
# few "Time LR" iteration counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "TestingFile" sheet (tabSelected, holds the query table)

# ---------------------------------------------------------------------
# 1. Snapshot the existing 12-column x 7-row table (header + 6 data rows)
#    before we start overwriting cells, since several columns move.
# ---------------------------------------------------------------------
$oldCols = 12
$oldRows = 7

$old = @{}
for ($r = 1; $r -le $oldRows; $r++) {
    for ($c = 1; $c -le $oldCols; $c++) {
        $old[[string]$r + "_" + [string]$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# Old layout (1-indexed columns):
#  1 Omega  2 Tolerance  3 Step Size Rule  4 Gap LR  5 Iterations
#  6 Converged?  7 Obj. Naive  8 Obj. LR  9 Gap  10 Time Naive
#  11 Time LR  12 Final Lambda

# ---------------------------------------------------------------------
# 2. Expand the ListObject (structured table) to the new 13-column range.
#    Resize() rewrites ref/autoFilter/tableColumns for us.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:M7"))

# ---------------------------------------------------------------------
# 3. New header row (1-indexed columns):
#  1 Omega  2 Tolerance  3 Step Size Rule  4 Iterations  5 Converged?
#  6 Gap LR  7 Gap Naive  8 Obj. Naive  9 Obj. LR  10 Gap  11 Time Naive
#  12 Time LR  13 Final Lambda
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Omega"
$ws.Cells.Item(1, 2).Value = "Tolerance"
$ws.Cells.Item(1, 3).Value = "Step Size Rule"
$ws.Cells.Item(1, 4).Value = "Iterations"
$ws.Cells.Item(1, 5).Value = "Converged?"
$ws.Cells.Item(1, 6).Value = "Gap LR"
$ws.Cells.Item(1, 7).Value = "Gap Naive"
$ws.Cells.Item(1, 8).Value = "Obj. Naive"
$ws.Cells.Item(1, 9).Value = "Obj. LR"
$ws.Cells.Item(1, 10).Value = "Gap"
$ws.Cells.Item(1, 11).Value = "Time Naive"
$ws.Cells.Item(1, 12).Value = "Time LR"
$ws.Cells.Item(1, 13).Value = "Final Lambda"

# ---------------------------------------------------------------------
# 4. Manual "Time LR" corrections captured from the run that regenerated
#    this scenario file (rows are keyed by the OLD row number, which is
#    unchanged since no rows were added/removed).
# ---------------------------------------------------------------------
$timeLrOverride = @{ 3 = 14; 5 = 17; 6 = 13 }

for ($r = 2; $r -le $oldRows; $r++) {
    $omega      = $old[[string]$r + "_1"]
    $tolerance  = $old[[string]$r + "_2"]
    $stepRule   = $old[[string]$r + "_3"]
    $gapLr      = $old[[string]$r + "_4"]
    $iterations = $old[[string]$r + "_5"]
    $converged  = $old[[string]$r + "_6"]
    $objNaive   = $old[[string]$r + "_7"]
    $objLr      = $old[[string]$r + "_8"]
    $gap        = $old[[string]$r + "_9"]
    $timeNaive  = $old[[string]$r + "_10"]
    $timeLr     = $old[[string]$r + "_11"]
    $finalLambda = $old[[string]$r + "_12"]

    if ($timeLrOverride.ContainsKey($r)) {
        $timeLr = $timeLrOverride[$r]
    }

    $ws.Cells.Item($r, 1).Value  = $omega
    $ws.Cells.Item($r, 2).Value  = $tolerance
    $ws.Cells.Item($r, 3).Value  = $stepRule
    $ws.Cells.Item($r, 4).Value  = $iterations
    $ws.Cells.Item($r, 5).Value  = $converged
    $ws.Cells.Item($r, 6).Value  = $gapLr
    $ws.Cells.Item($r, 7).Value  = 0            # new "Gap Naive" column
    $ws.Cells.Item($r, 8).Value  = $objNaive
    $ws.Cells.Item($r, 9).Value  = $objLr
    $ws.Cells.Item($r, 10).Value = $gap
    $ws.Cells.Item($r, 11).Value = $timeNaive
    $ws.Cells.Item($r, 12).Value = $timeLr
    $ws.Cells.Item($r, 13).Value = $finalLambda
}

# ---------------------------------------------------------------------
# 5. The hidden "ExterneDaten_1" (External Data) defined name tracks the
#    query table's range and must grow from L7 to M7 along with the table.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    $nm.RefersTo = "=TestingFile!`$A`$1:`$M`$7"
}
